# Update "想去人数" (F column) counts on both the "展览" sheet and the
# "全部类型" aggregate sheet, as published by the latest data refresh.

$wb = $excel.ActiveWorkbook

$exhibitionSheet = $wb.Worksheets.Item("展览")
$allTypesSheet   = $wb.Worksheets.Item("全部类型")

# Row -> new value, keyed by sheet (rows differ by one offset for rows
# at/after row 8 on the "展览" sheet because "全部类型" has one extra
# row of content mixed in above that point).
$exhibitionUpdates = @{
    2  = 640
    3  = 2231
    5  = 13416
    8  = 525
    10 = 1199
    13 = 14521
    22 = 47
    27 = 5561
    28 = 944
    31 = 34
    33 = 148
}

$allTypesUpdates = @{
    2  = 640
    3  = 2231
    5  = 13416
    9  = 525
    11 = 1199
    14 = 14521
    23 = 47
    28 = 5561
    29 = 944
    32 = 34
    34 = 148
}

foreach ($row in $exhibitionUpdates.Keys) {
    $exhibitionSheet.Range("F$row").Value = $exhibitionUpdates[$row]
}

foreach ($row in $allTypesUpdates.Keys) {
    $allTypesSheet.Range("F$row").Value = $allTypesUpdates[$row]
}
